# The submission letter's dateline reads "July 25, 2025". Update it to
# "September 11, 2025" (only the "July 25" portion changes; the trailing
# ", 2025" stays as-is).
$d = $word.ActiveDocument

$find = $d.Content.Find
$find.Execute("July 25", $true, $false, $false, $false, $false, $true, 1, $false, "September 11", 2)
